$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 values (A1:N1 changed, O1:Q1 are new)
$ws.Range("A1").Value = 5
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 27
$ws.Range("E1").Value = 20
$ws.Range("F1").Value = 18
$ws.Range("G1").Value = 32
$ws.Range("H1").Value = 32
$ws.Range("I1").Value = 16
$ws.Range("J1").Value = 14
$ws.Range("K1").Value = 3
$ws.Range("L1").Value = 26
$ws.Range("M1").Value = 0.056999999999999995
$ws.Range("N1").Value = 0.012
$ws.Range("O1").Value = 0.053000000000000005
$ws.Range("P1").Value = 0.074999999999999997
$ws.Range("Q1").Value = 0.050000000000000003

# Update column widths:
#  - G (7) widens from the narrow width to the regular width
#  - K (11) narrows from the wide width down to the narrow width
#  - L (12) narrows from the wide width down to the regular width
#  - new columns O:Q (15-17) get the wide data-column widths
# (ColumnWidth is expressed in characters and Excel snaps it to its pixel
# grid, so these inputs are chosen to land on the nearest grid value to the
# target column widths.)
$ws.Columns.Item(7).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(11).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(12).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(15).ColumnWidth = 4.833333333333333
$ws.Columns.Item(16).ColumnWidth = 4.833333333333333
$ws.Columns.Item(17).ColumnWidth = 3.8333333333333335
